{"js": "// The document has a bulleted list of \"games referenced during design\",\n// e.g. \u7089\u77f3\u4f20\u8bf4 / \u6697\u9ed1\u7834\u574f\u795e3 / \u690d\u7269\u5927\u6218\u50f5\u5c38 / \u90e8\u843d\u51b2\u7a81 / \u5e1d\u56fd\u65f6\u4ee3 /\n// \u575a\u5b88\u9635\u57302 / \u63a8\u7bb1\u5b50 / \u76df\u519b\u6562\u6b7b\u961f (in that exact order, back to back).\n// The edit appends a short \"\uff1a<note>\" explanation to each of those eight\n// bullet paragraphs (commit: \"split into 2 class for attacker and\n// defenser\").\n//\n// Map of exact paragraph text (as it exists today) -> ordered list of text\n// chunks to append at the end of the paragraph (i.e. before the paragraph\n// mark). Keeping the chunks as separate insertText() calls mirrors the\n// separate <w:r> runs seen in the authoritative diff as closely as the\n// Word JS API allows.\nconst additions = [\n  [\"\u7089\u77f3\u4f20\u8bf4\", [\"\uff1a\", \"\u6570\u5b57\u663e\u793a\", \"\u8840\u91cf\u548c\u5355\u4f4d\", \"\uff0c\", \"\u653b\u51fb\u673a\u5236\u3002\"]],\n  [\"\u6697\u9ed1\u7834\u574f\u795e3\", [\"\uff1a\u8303\u56f4\u653b\u51fb\uff0c\", \"PVE\", \"\u7684\u4e50\u8da3\u70b9\u3002\"]],\n  [\"\u690d\u7269\u5927\u6218\u50f5\u5c38\", [\"\uff1a\", \"\u683c\u5b50\", \"\u5373\u65f6\u5236\", \"\uff0c\", \"\u7c7b\", \"\u5854\u9632\", \"\u3002\"]],\n  [\"\u90e8\u843d\u51b2\u7a81\", [\"\uff1a\", \"\u5373\u65f6\u5236\", \"+\", \"\u9006\u5411\u5854\u9632\u3002\"]],\n  [\"\u5e1d\u56fd\u65f6\u4ee3\", [\"\uff1a\", \"\u5355\u4f4d\", \"AI\", \"\u3002\"]],\n  [\"\u575a\u5b88\u9635\u57302\", [\"\uff1a\", \"\u7c7b\u5854\u9632\"]],\n  [\"\u63a8\u7bb1\u5b50\", [\"\uff1a\", \"\u79fb\u52a8\", \"\u5355\u4f4d\", \"\u3002\"]],\n  [\"\u76df\u519b\u6562\u6b7b\u961f\", [\"\uff1a\", \"\u5355\u4f4d\", \"AI\", \"\uff0c\", \"\u5355\u4f4d\", \"\u76d1\u89c6\u8303\u56f4\"]],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The eight bullets we care about all live next to each other (they're the\n// \"\u8bbe\u8ba1\u8fc7\u7a0b\u4e2d\u53c2\u8003\u8fc7\u7684\u6e38\u620f\" sub-list); match by exact trimmed text so we do\n// not also catch unrelated paragraphs elsewhere in the doc that merely\n// contain the same game name as a substring (e.g. \"\u63a8\u7bb1\u5b50\"/\"\u76df\u519b\u6562\u6b7b\u961f\"\n// each also show up quoted inside a different bullet earlier in the file).\nlet additionIndex = 0;\nfor (let i = 0; i < paragraphs.items.length && additionIndex < additions.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text.trim();\n  const [target, chunks] = additions[additionIndex];\n  if (text === target) {\n    for (const chunk of chunks) {\n      para.insertText(chunk, \"End\");\n    }\n    additionIndex++;\n  }\n}\n\nawait context.sync();\n\nif (additionIndex !== additions.length) {\n  throw new Error(\n    `Only matched ${additionIndex} of ${additions.length} target paragraphs`\n  );\n}\n", "ps1": "# The document has a bulleted list of \"games referenced during design\",\n# e.g. \u7089\u77f3\u4f20\u8bf4 / \u6697\u9ed1\u7834\u574f\u795e3 / \u690d\u7269\u5927\u6218\u50f5\u5c38 / \u90e8\u843d\u51b2\u7a81 / \u5e1d\u56fd\u65f6\u4ee3 /\n# \u575a\u5b88\u9635\u57302 / \u63a8\u7bb1\u5b50 / \u76df\u519b\u6562\u6b7b\u961f (in that exact order, back to back).\n# The edit appends a short \"\uff1a<note>\" explanation to each of those eight\n# bullet paragraphs (commit: \"split into 2 class for attacker and\n# defenser\").\n#\n# Build an ordered list of (exact current paragraph text, chunks to append)\n# pairs. Appending the chunks one at a time (instead of one big\n# concatenated string) mirrors the separate <w:r> runs seen in the\n# authoritative diff as closely as the Word object model allows.\n$targets = @(\n    @{ Name = \"\u7089\u77f3\u4f20\u8bf4\"; Parts = @(\"\uff1a\", \"\u6570\u5b57\u663e\u793a\", \"\u8840\u91cf\u548c\u5355\u4f4d\", \"\uff0c\", \"\u653b\u51fb\u673a\u5236\u3002\") },\n    @{ Name = \"\u6697\u9ed1\u7834\u574f\u795e3\"; Parts = @(\"\uff1a\u8303\u56f4\u653b\u51fb\uff0c\", \"PVE\", \"\u7684\u4e50\u8da3\u70b9\u3002\") },\n    @{ Name = \"\u690d\u7269\u5927\u6218\u50f5\u5c38\"; Parts = @(\"\uff1a\", \"\u683c\u5b50\", \"\u5373\u65f6\u5236\", \"\uff0c\", \"\u7c7b\", \"\u5854\u9632\", \"\u3002\") },\n    @{ Name = \"\u90e8\u843d\u51b2\u7a81\"; Parts = @(\"\uff1a\", \"\u5373\u65f6\u5236\", \"+\", \"\u9006\u5411\u5854\u9632\u3002\") },\n    @{ Name = \"\u5e1d\u56fd\u65f6\u4ee3\"; Parts = @(\"\uff1a\", \"\u5355\u4f4d\", \"AI\", \"\u3002\") },\n    @{ Name = \"\u575a\u5b88\u9635\u57302\"; Parts = @(\"\uff1a\", \"\u7c7b\u5854\u9632\") },\n    @{ Name = \"\u63a8\u7bb1\u5b50\"; Parts = @(\"\uff1a\", \"\u79fb\u52a8\", \"\u5355\u4f4d\", \"\u3002\") },\n    @{ Name = \"\u76df\u519b\u6562\u6b7b\u961f\"; Parts = @(\"\uff1a\", \"\u5355\u4f4d\", \"AI\", \"\uff0c\", \"\u5355\u4f4d\", \"\u76d1\u89c6\u8303\u56f4\") }\n)\n\n$d = $word.ActiveDocument\n$matchIndex = 0\n$paraCount = $d.Paragraphs.Count\n\nfor ($i = 1; $i -le $paraCount -and $matchIndex -lt $targets.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.Trim()\n    $target = $targets[$matchIndex]\n\n    if ($text -eq $target.Name) {\n        $r = $p.Range\n        # Exclude the trailing paragraph mark so text lands at the end of\n        # the visible paragraph text, not at the start of the next one.\n        $r.MoveEnd(1, -1)\n        foreach ($part in $target.Parts) {\n            $r.InsertAfter($part)\n            $r.Collapse(0)\n        }\n        $matchIndex = $matchIndex + 1\n    }\n}\n\nWrite-Host \"Matched\" $matchIndex \"of\" $targets.Count \"target paragraphs\"\n"}
